# Weekly update: insert a new data row for "Pepino ensalada" (Macroferia
# Regional de Talca) at row 341, pushing the existing rows 341-435 down to
# 342-436 (the sheet grows from A1:R435 to A1:R436).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 341; this shifts rows
# 341..435 down to 342..436 and preserves their values/formatting.
$ws.Rows(341).Insert()

# Populate the newly inserted row 341 with this week's new record.
$ws.Range("A341").Value = 5
$ws.Range("B341").Value = "Macroferia Regional de Talca"
$ws.Range("C341").Value = "Maule"
$ws.Range("D341").Value = 44736
$ws.Range("E341").Value = 7
$ws.Range("F341").Value = 100112043
$ws.Range("G341").Value = "Pepino ensalada"
$ws.Range("H341").Value = "Sin especificar"
$ws.Range("I341").Value = "Primera"
$ws.Range("J341").Value = 300
$ws.Range("K341").Value = 19000
$ws.Range("L341").Value = 19000
$ws.Range("M341").Value = 19000
$ws.Range("N341").Value = "`$/caja 60 unidades"
$ws.Range("O341").Value = "Región de Arica y Parinacota"
$ws.Range("P341").Value = 317
$ws.Range("Q341").Value = 60
$ws.Range("R341").Value = "Hortaliza"
